# job_history.xlsx edit:
#   "run 10% subset with less columns, results are also not good"
#
# - Row 29 (previously a blank spacer row) is filled in with a new
#   "10% subset, fewer columns" result row, reusing the same category
#   values as row 27 (ukb51139_subset.csv / all / no events / > 140/80 /
#   zscore / median / none) but with the "2801 x 145" subset dimensions,
#   Vars Used = 51, new 1/0 and 2/0 overlap results, # cPCA = 17 and
#   Alpha = 2.66.
# - Rows 27, 28, 30 and 31 have their K/N/O (and I for the blank rows)
#   result cells recolored from the "theme" black font to an explicit
#   black font -- the "results are also not good" visual flag -- which
#   is the same formatting already used on row 29's sibling rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: clone row 27's formatting (borders/number formats/fonts
# for every column), then overwrite with the new row's own values. ---
$ws.Range("A27:O27").Copy() | Out-Null
$ws.Range("A29:O29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(29, 1).Value = "ukb51139_subset.csv"
$ws.Cells.Item(29, 2).Value = "2801 x 145"
$ws.Cells.Item(29, 3).Value = "all"
$ws.Cells.Item(29, 4).Value = "no events"
$ws.Cells.Item(29, 5).Value = "> 140/80"
$ws.Cells.Item(29, 6).Value = "zscore"
$ws.Cells.Item(29, 7).Value = "median"
$ws.Cells.Item(29, 8).Value = "none"
$ws.Cells.Item(29, 9).Value = 50
$ws.Cells.Item(29, 11).Value = 51
$ws.Cells.Item(29, 12).Value = "93.1 & 82.5"
$ws.Cells.Item(29, 13).Value = "52.1 & 51.4"
$ws.Cells.Item(29, 14).Value = 17
$ws.Cells.Item(29, 15).Value = 2.66

# --- Rows 27 & 28: recolor the Vars Used / # cPCA / Alpha cells to the
# explicit-black "not good" font. ---
$ws.Cells.Item(27, 11).Font.Color = 0
$ws.Cells.Item(27, 14).Font.Color = 0
$ws.Cells.Item(27, 15).Font.Color = 0

$ws.Cells.Item(28, 11).Font.Color = 0
$ws.Cells.Item(28, 14).Font.Color = 0
$ws.Cells.Item(28, 15).Font.Color = 0

# Row 28 also grew slightly taller to match the other data rows.
$ws.Rows.Item(28).RowHeight = 19.5

# --- Rows 30 & 31: same recolor, including the (still blank) I column. ---
$ws.Cells.Item(30, 9).Font.Color = 0
$ws.Cells.Item(30, 11).Font.Color = 0
$ws.Cells.Item(30, 14).Font.Color = 0
$ws.Cells.Item(30, 15).Font.Color = 0

$ws.Cells.Item(31, 9).Font.Color = 0
$ws.Cells.Item(31, 11).Font.Color = 0
$ws.Cells.Item(31, 14).Font.Color = 0
$ws.Cells.Item(31, 15).Font.Color = 0
